$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 19: "Linear Gesture" test case (RE2.2 Generic DSL)
# ---------------------------------------------------------------------
$srcRange = $ws.Range("A18:K18")
$srcRange.Copy()
$dstRange19 = $ws.Range("A19:K19")
$dstRange19.PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 192

$steps19 = @'
wait(3);
validate1;
link_Click(EMML11_test_link);
wait(2);
validate2;
link_Click(gesture11_test_link);
DrawGesture(linear_default);
validate3;
'@

$validate19 = @'
validate1
{
validate_PageTitle=Pocket Browser Tests
};
validate2
{
validate_PageTitle=EMML1.1 Test Index Page
};
validate3
{
validate_Result=Gesture detected for 1th time
validate_Result=swipe
};
'@

$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = 1
$ws.Cells.Item(19,4).Value = "A"
$ws.Cells.Item(19,5).Value = "Linear Gesture"
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = $steps19
$ws.Cells.Item(19,8).Value = $validate19

# ---------------------------------------------------------------------
# Row 20: "Hold Gesture" test case (RE2.2 Generic DSL)
# ---------------------------------------------------------------------
$srcRange.Copy()
$dstRange20 = $ws.Range("A20:K20")
$dstRange20.PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 192

$steps20 = @'
wait(3);
validate1;
link_Click(EMML11_test_link);
wait(2);
validate2;
link_Click(gesture11_test_link);
DrawGesture(hold,100,100,6000);
validate3;
'@

$validate20 = @'
validate1
{
validate_PageTitle=Pocket Browser Tests
};
validate2
{
validate_PageTitle=EMML1.1 Test Index Page
};
validate3
{
validate_Result=Gesture detected for 3th time
validate_Result=press
};
'@

$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = 1
$ws.Cells.Item(20,4).Value = "A"
$ws.Cells.Item(20,5).Value = "Hold Gesture"
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = $steps20
$ws.Cells.Item(20,8).Value = $validate20
